# Project Euler 005 - add a third solving technique (SEQUENCE+REDUCE one-liner),
# renumber/relabel the existing sections, and duplicate the REDUCE example into
# the "array formulas" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 0) Make room: the old "4) Use of spreadsheets capabilities" block
#    (rows 25-46) needs to move down by 5 rows to 30-51, to leave space
#    for a new "array formulas" demo block (rows 25-29).
# ------------------------------------------------------------------
$ws.Rows("25:29").Insert()

# ------------------------------------------------------------------
# 1) Capture the "section title" look (bold/colored, no quote-prefix)
#    that currently lives on B18, and stamp it onto B23 BEFORE we
#    overwrite B18's own text/format in step 2. (B23 keeps referencing
#    the same slot conceptually, it just becomes the new
#    "array formulas" header once B18 becomes the SEQUENCE+REDUCE one.)
# ------------------------------------------------------------------
$ws.Range("B18").Copy()
$ws.Range("B23").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B23").Value2 = "4) One-liner based on array formulas"

# ------------------------------------------------------------------
# 2) New section 3 header + its surrounding spacer rows, styled like
#    the "2) tail-call recursion" header (B13:D13 / B14:D14 / B12:D12).
# ------------------------------------------------------------------
$ws.Range("B12:D12").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B13:D13").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B18").Value2 = "3) One-liner converting recursion into SEQUENCE+REDUCE"

$ws.Range("B14:D14").Copy()
$ws.Range("B19:D19").PasteSpecial(-4122)  # xlPasteFormats

# ------------------------------------------------------------------
# 3) Rewrite the REDUCE one-liner (row 21) with renamed LAMBDA params
#    (a,v -> ACC,N) and Excel's re-pretty-printed spacing.
# ------------------------------------------------------------------
$ws.Range("B21").Formula = " =REDUCE( 1, SEQUENCE(B20), LAMBDA(ACC,N, LCM(ACC,N)))"
$ws.Range("D21").Formula = " =REDUCE( 1, SEQUENCE(D20), LAMBDA(ACC,N, LCM(ACC,N)))"

$ws.Range("B12:D12").Copy()
$ws.Range("B22:D22").PasteSpecial(-4122)  # xlPasteFormats

# ------------------------------------------------------------------
# 4) New "array formulas" demo block (rows 25-26): same pattern as the
#    10/20-seeded REDUCE example above, duplicated verbatim.
# ------------------------------------------------------------------
$ws.Range("B20:D20").Copy()
$ws.Range("B25:D25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B25").Value2 = 10
$ws.Range("D25").Value2 = 20

$ws.Range("B21:E21").Copy()
$ws.Range("B26:E26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B26").Formula = "=REDUCE(1, SEQUENCE(B25), LAMBDA(a,v, LCM(a,v)))"
$ws.Range("D26").Formula = "=REDUCE(1, SEQUENCE(D25), LAMBDA(a,v, LCM(a,v)))"
$ws.Range("E26").Value2 = $ws.Range("E21").Value2

# ------------------------------------------------------------------
# 5) New section title for the (renumbered) "Use of spreadsheets
#    capabilities" block, styled like "1) One-liner based on recursion".
# ------------------------------------------------------------------
$ws.Range("B8").Copy()
$ws.Range("B28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B28").Value2 = "5) Use of spreadsheets capabilities"

# ------------------------------------------------------------------
# 6) Cosmetic: selection moved to G3 in the saved file.
# ------------------------------------------------------------------
$ws.Range("G3").Select()
